# Weekly update: insert one new price record at the top of the
# "Femacal de La Calera - Mango" data block (row 610), pushing the
# existing rows 610-727 down to 611-728.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 610 (shifts 610:727 -> 611:728).
$ws.Rows.Item(610).Insert()

# Populate the new row 610 with this week's record.
$ws.Cells.Item(610, 1).Value  = 3                                   # A Mercado ID
$ws.Cells.Item(610, 2).Value  = "Femacal de La Calera"               # B Mercado
$ws.Cells.Item(610, 3).Value  = "Coquimbo"                           # C Región
$ws.Cells.Item(610, 4).Value  = 45209                                # D Fecha
$ws.Cells.Item(610, 5).Value  = 5                                    # E Codreg
$ws.Cells.Item(610, 6).Value  = "Fruta"                              # F Tipo
$ws.Cells.Item(610, 7).Value  = 100108                               # G Producto ID
$ws.Cells.Item(610, 8).Value  = "Tropicales y subtropicales"         # H Producto
$ws.Cells.Item(610, 9).Value  = 100108002                            # I Categoría ID
$ws.Cells.Item(610, 10).Value = "Mango"                              # J Categoría
$ws.Cells.Item(610, 11).Value = "Sin especificar"                    # K Variedad
$ws.Cells.Item(610, 12).Value = "Primera"                            # L Calidad
$ws.Cells.Item(610, 13).Value = 228                                  # M Volumen
$ws.Cells.Item(610, 14).Value = 9000                                 # N Precio mínimo
$ws.Cells.Item(610, 15).Value = 9000                                 # O Precio máximo
$ws.Cells.Item(610, 16).Value = 9000                                 # P Precio promedio ponderado
$ws.Cells.Item(610, 17).Value = "$/bandeja 4 kilos"                  # Q Unidad de comercialización
$ws.Cells.Item(610, 18).Value = "Brasil"                             # R Origen
$ws.Cells.Item(610, 19).Value = 2250                                 # S Precio $/Kg
$ws.Cells.Item(610, 20).Value = 4                                    # T Kg / unidad

# Keep the date cell's date/time number format consistent with the rest
# of column D.
$ws.Cells.Item(610, 4).NumberFormat = $ws.Cells.Item(611, 4).NumberFormat
